$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C:E to D:F)
$ws.Columns("C").Insert()

# New header for inserted column
$ws.Range("C1").Value = "Cek role"

# Select C2 to match final selection state
$ws.Range("C2").Select()
